# The edit: a requirements document gains a new run ("lll", Arial 8pt,
# pt-BR) immediately before the pre-existing run (two spaces, same
# formatting) inside the very last paragraph of the document body.
#
# A plain Range.InsertBefore()/Range.Text=/Find.Execute() insertion here
# would get silently coalesced into the neighbouring run (same rPr), so
# the new text would end up glued onto the existing run instead of
# living in its own <w:r>. Range.InsertXML() inserts pre-built run XML
# verbatim and keeps the run boundary intact, which is what the diff
# calls for.

$d = $word.ActiveDocument

$p = $d.Paragraphs.Last
$r = $p.Range

if ($p.Range.Text -ne "  `r") {
    throw "Unexpected target paragraph text: '$($p.Range.Text)'"
}

# Collapsed range sitting right before the existing run, i.e. at the
# very start of the paragraph.
$insertionPoint = $d.Range($r.Start, $r.Start)

$xmlFragment = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                <w:sz w:val="16"/>
                <w:lang w:val="pt-BR"/>
              </w:rPr>
              <w:t>lll</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

[void]$insertionPoint.InsertXML($xmlFragment)

Write-Output "Inserted 'lll' run before the trailing paragraph's existing run."
